# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
#
# The "Date" column (BF) on Sheet1 held the source-filename-derived string
# "6-16-2007-08" for every row (rows 2-31); correct it to the actual game
# date "2008-06-16" (ISO yyyy-mm-dd), keeping the value as literal text
# (not an auto-converted date serial number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("BF2:BF31")
$rowCount = $target.Rows.Count

# Stage the corrected text in a scratch range that is explicitly formatted
# as Text, so Excel's smart-entry parser doesn't turn the ISO-looking
# string into a date when it is (re)written. Then copy/paste-special just
# the values into the real destination, which carries across the literal
# text without disturbing the destination cells' existing (default)
# number format/style.
$scratch = $ws.Range("ZZ1").Resize($rowCount, 1)
$scratch.NumberFormat = "@"
for ($i = 1; $i -le $rowCount; $i++) {
    $scratch.Cells.Item($i).Value = "2008-06-16"
}

$scratch.Copy()
$target.PasteSpecial(-4163)  # xlPasteValues

$scratch.Clear()
$excel.CutCopyMode = $false
